# Automatische test-sync: 2025-06-26 23:31:50
#
# Adds the new "Testmail #10" row to the Logs sheet, rolls the matching
# tally into the Dashboard sheet, and extends the bar chart's category /
# value series references to cover the new Dashboard row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 42 (Testmail #10: "Zijn er vacatures?")
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A42").Value = "Zijn er vacatures?"
$logs.Range("B42").Value = "mailmind.test@zohomail.eu"
$logs.Range("C42").Value = "Testmail #10: Zijn er vacatures?"
$logs.Range("D42").Value = "Sollicitatie / Vacature"
$logs.Range("E42").Value = "Beste,`nDank u voor uw interesse in onze organisatie. Op dit moment hebben we geen openstaande vacatures, maar we moedigen u aan om regelmatig onze website te bezoeken voor eventuele toekomstige mogelijkheden. Mocht u nog vragen hebben of meer informatie wensen, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Range("F42").Value = "2025-06-26 23:30:58"
$logs.Range("G42").Value = "Ja"
$logs.Range("H42").Value = "Nee"
$logs.Range("I42").Value = "Ja"

# E42 contains embedded line breaks, which makes the host auto-pin an
# explicit (wrapped) row height. Re-run AutoFit so the row reverts to the
# sheet's implicit default height, matching every other row in the sheet.
$logs.Rows.Item(42).AutoFit()

# The sheet's used-range dimension grows automatically, but the existing
# conditional-formatting rules keep their old sqref (row 2-41) unless we
# explicitly re-apply them to the new row 2-42 range.
$logs.Range("D2:D41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D42"))
$logs.Range("G2:G41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G42"))
$logs.Range("H2:H41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H42"))
$logs.Range("I2:I41").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I42"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append row 9 (new "Sollicitatie / Vacature" tally)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A9").Value = "Sollicitatie / Vacature"
$dash.Range("B9").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the bar chart's category/value series to row 9
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.XValues = "='Dashboard'!`$A`$2:`$A`$9"
$ser.Values = "='Dashboard'!`$B`$2:`$B`$9"
